$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven login automation: read the usernames/passwords already in
# columns A/B and write the authentication result ("Correct"/"Incorrect")
# into column C for each row. Row 2's password is treated as the correct
# one that subsequent attempts are checked against.
$expected = $ws.Range("B2").Value2

$ws.Range("C2").Value = "Correct"

if ($ws.Range("B3").Value2 -eq $expected) {
    $ws.Range("C3").Value = "Correct"
} else {
    $ws.Range("C3").Value = "Incorrect"
}

if ($ws.Range("B4").Value2 -eq $expected) {
    $ws.Range("C4").Value = "Correct"
} else {
    $ws.Range("C4").Value = "Incorrect"
}

if ($ws.Range("B5").Value2 -eq $expected) {
    $ws.Range("C5").Value = "Correct"
} else {
    $ws.Range("C5").Value = "Incorrect"
}
